# Update config for light tomo (odin -> ymir topics/sources)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: rotation_angle (NXsample/sample) - motion topic + PV source
$ws.Range("C2").Value = "ymir_motion"
$ws.Range("D2").Value = "SES-SCAN:MC-MCU-001:m4.RBV"

# Row 3: NXmonitor/control - camera topic + new source
$ws.Range("C3").Value = "ymir_camera"
$ws.Range("D3").Value = "other_source"

# Row 5: NXdetector/detector image data - camera topic, new source, new array size
$ws.Range("C5").Value = "ymir_camera"
$ws.Range("D5").Value = "some_source"
$ws.Range("H5").Value = "650, 650"

# Row 6: NXdetector/detector image_key - camera topic
$ws.Range("C6").Value = "ymir_camera"

# Column width adjustments (no longer best-fit, explicit custom widths)
# (values compensate for this engine's pixel-rounding of ColumnWidth so the
# stored OOXML width lands on the target: col C -> 12, col D -> 26.33203125)
$ws.Columns.Item(3).ColumnWidth = 11.17
$ws.Columns.Item(4).ColumnWidth = 25.5

# Update selection to match the new active cell
$ws.Range("E6").Select()
